# GameGetter Config.xlsx: add EpicCredential / SteamCredential settings
# rows to the "Settings" sheet (ahead of the existing Steam_Credential /
# Steam_Failed_Credential rows), and update the active sheet/selection to
# reflect the editor having been working on the Settings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Insert two new rows above the existing "Steam_Credential" row (row 20),
# pushing the Steam_Credential / Steam_Failed_Credential rows down to 22/23.
$ws.Rows("19:20").Insert()

# New row 19: EpicCredential setting (Name / Value / Description).
$ws.Range("A19").RowHeight = 14.25
$ws.Range("A19").Value = "EpicCredential"
$ws.Range("B19").Value = "EpicCredential"
$ws.Range("C19").Value = "Name for credential to sign into Epic Games"

# Row 20 stays blank (spacer row), matching the sheet's existing layout
# convention of a blank separator row between each Name/Value/Description
# entry.
$ws.Range("A20").RowHeight = 14.25

# New row 21: SteamCredential setting (only Name / Value are populated).
$ws.Range("A21").Value = "SteamCredential"
$ws.Range("B21").Value = "SteamCredential"

# Bring focus to the Settings sheet (was Assets) and leave the selection
# where the author left off editing.
$ws.Activate() | Out-Null
$ws.Range("B26").Select() | Out-Null
